# Refresh market-price derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# for a batch of leve rows across multiple sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H87").Value = 55397
$ws.Range("J87").Value = 58862.668
$ws.Range("L87").Value = 58862.668
$ws.Range("N87").Value = -61358.668

$ws.Range("H90").Value = 55397
$ws.Range("J90").Value = 58862.668
$ws.Range("L90").Value = 176588.004
$ws.Range("N90").Value = -189068.004

$ws.Range("H96").Value = 727.1111
$ws.Range("I96").Value = 727.1111
$ws.Range("K96").Value = 2181.3333
$ws.Range("M96").Value = -808.3332999999998

$ws.Range("H118").Value = 319.07693
$ws.Range("I118").Value = 304
$ws.Range("K118").Value = 912
$ws.Range("M118").Value = 745

$ws.Range("H127").Value = 5548.6665
$ws.Range("I127").Value = 4860.4
$ws.Range("K127").Value = 14581.2
$ws.Range("M127").Value = -9621.199999999999

$ws.Range("H131").Value = 3637.5386
$ws.Range("I131").Value = 2274
$ws.Range("J131").Value = 20000
$ws.Range("K131").Value = 6822
$ws.Range("L131").Value = 60000
$ws.Range("M131").Value = -1782
$ws.Range("N131").Value = -70080

$ws.Range("H137").Value = 3856.3333
$ws.Range("I137").Value = 3992.0344
$ws.Range("J137").Value = 2872.5
$ws.Range("K137").Value = 11976.1032
$ws.Range("L137").Value = 8617.5
$ws.Range("M137").Value = -9426.1032
$ws.Range("N137").Value = -13717.5

$ws.Range("H138").Value = 2390.8936
$ws.Range("J138").Value = 3474.4
$ws.Range("L138").Value = 10423.2
$ws.Range("N138").Value = -20703.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 21733.334
$ws.Range("J55").Value = 31100
$ws.Range("L55").Value = 31100
$ws.Range("N55").Value = -31730

$ws.Range("H61").Value = 5872.2383
$ws.Range("I61").Value = 5986.025
$ws.Range("K61").Value = 5986.025
$ws.Range("M61").Value = -5774.025

$ws.Range("H74").Value = 4256.864
$ws.Range("J74").Value = 8200
$ws.Range("L74").Value = 8200
$ws.Range("N74").Value = -9948

$ws.Range("H77").Value = 4256.864
$ws.Range("J77").Value = 8200
$ws.Range("L77").Value = 41000
$ws.Range("N77").Value = -49736

$ws.Range("H97").Value = 837.44183
$ws.Range("I97").Value = 742.9211
$ws.Range("K97").Value = 742.9211
$ws.Range("M97").Value = -246.9211

$ws.Range("H122").Value = 2392.3333
$ws.Range("I122").Value = 2392.3333
$ws.Range("K122").Value = 7176.999899999999
$ws.Range("M122").Value = -4726.999899999999

$ws.Range("H136").Value = 5872.2383
$ws.Range("I136").Value = 5986.025
$ws.Range("K136").Value = 17958.075
$ws.Range("M136").Value = -15408.075

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 32500
$ws.Range("J35").Value = 32500
$ws.Range("L35").Value = 32500
$ws.Range("N35").Value = -33120

$ws.Range("H82").Value = 15889.25
$ws.Range("J82").Value = 44300
$ws.Range("L82").Value = 44300
$ws.Range("N82").Value = -45066

$ws.Range("H85").Value = 15889.25
$ws.Range("J85").Value = 44300
$ws.Range("L85").Value = 44300
$ws.Range("N85").Value = -46952

$ws.Range("H94").Value = 1274.5625
$ws.Range("I94").Value = 499.83334
$ws.Range("J94").Value = 1739.4
$ws.Range("K94").Value = 499.83334
$ws.Range("L94").Value = 1739.4
$ws.Range("M94").Value = -48.83334000000002
$ws.Range("N94").Value = -2641.4

$ws.Range("H134").Value = 6161.391
$ws.Range("I134").Value = 5997.8203
$ws.Range("K134").Value = 17993.4609
$ws.Range("M134").Value = -15458.4609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4899.9287
$ws.Range("I31").Value = 3265.9167
$ws.Range("J31").Value = 6125.4375
$ws.Range("K31").Value = 3265.9167
$ws.Range("L31").Value = 6125.4375
$ws.Range("M31").Value = -2970.9167
$ws.Range("N31").Value = -6715.4375

$ws.Range("H34").Value = 4899.9287
$ws.Range("I34").Value = 3265.9167
$ws.Range("J34").Value = 6125.4375
$ws.Range("K34").Value = 3265.9167
$ws.Range("L34").Value = 6125.4375
$ws.Range("M34").Value = -3063.9167
$ws.Range("N34").Value = -6529.4375

$ws.Range("H134").Value = 5116.0557
$ws.Range("I134").Value = 3322.1
$ws.Range("J134").Value = 7358.5
$ws.Range("K134").Value = 9966.299999999999
$ws.Range("L134").Value = 22075.5
$ws.Range("M134").Value = -7431.299999999999
$ws.Range("N134").Value = -27145.5

$ws.Range("H141").Value = 195841.42
$ws.Range("J141").Value = 195841.42
$ws.Range("L141").Value = 195841.42
$ws.Range("N141").Value = -206201.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1551.7142
$ws.Range("I39").Value = 427
$ws.Range("J39").Value = 1739.1666
$ws.Range("K39").Value = 1281
$ws.Range("L39").Value = 5217.4998
$ws.Range("M39").Value = -987
$ws.Range("N39").Value = -5805.4998

$ws.Range("H87").Value = 10428.429
$ws.Range("I87").Value = 7999.8335
$ws.Range("J87").Value = 25000
$ws.Range("K87").Value = 23999.5005
$ws.Range("L87").Value = 75000
$ws.Range("M87").Value = -22751.5005
$ws.Range("N87").Value = -77496

$ws.Range("H90").Value = 10428.429
$ws.Range("I90").Value = 7999.8335
$ws.Range("J90").Value = 25000
$ws.Range("K90").Value = 71998.5015
$ws.Range("L90").Value = 225000
$ws.Range("M90").Value = -65758.5015
$ws.Range("N90").Value = -237480

$ws.Range("H132").Value = 2474.5
$ws.Range("J132").Value = 2555.5557
$ws.Range("L132").Value = 23000.0013
$ws.Range("N132").Value = -28060.0013

$ws.Range("H134").Value = 1107.7894
$ws.Range("I134").Value = 1107.7894
$ws.Range("K134").Value = 3323.3682
$ws.Range("M134").Value = 1746.6318

$ws.Range("H137").Value = 5005
$ws.Range("J137").Value = 6426.625
$ws.Range("L137").Value = 19279.875
$ws.Range("N137").Value = -29479.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 13999
$ws.Range("I29").Value = 13999
$ws.Range("K29").Value = 13999
$ws.Range("M29").Value = -13709

$ws.Range("H97").Value = 740.6875
$ws.Range("J97").Value = 820.5
$ws.Range("L97").Value = 820.5
$ws.Range("N97").Value = -1812.5

$ws.Range("H113").Value = 2589.8667
$ws.Range("I113").Value = 2664.2856
$ws.Range("K113").Value = 2664.2856
$ws.Range("M113").Value = -494.2856000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 7000
$ws.Range("I45").Value = 6000
$ws.Range("K45").Value = 6000
$ws.Range("M45").Value = -5593

$ws.Range("H61").Value = 41948.2
$ws.Range("I61").Value = 45500.26
$ws.Range("J61").Value = 1099.5
$ws.Range("K61").Value = 45500.26
$ws.Range("L61").Value = 1099.5
$ws.Range("M61").Value = -45298.26
$ws.Range("N61").Value = -1503.5

$ws.Range("H68").Value = 5175.5
$ws.Range("I68").Value = 3234
$ws.Range("J68").Value = 11000
$ws.Range("K68").Value = 3234
$ws.Range("L68").Value = 11000
$ws.Range("M68").Value = -2485
$ws.Range("N68").Value = -12498

$ws.Range("H71").Value = 5175.5
$ws.Range("I71").Value = 3234
$ws.Range("J71").Value = 11000
$ws.Range("K71").Value = 16170
$ws.Range("L71").Value = 55000
$ws.Range("M71").Value = -12426
$ws.Range("N71").Value = -62488

$ws.Range("H113").Value = 41948.2
$ws.Range("I113").Value = 45500.26
$ws.Range("J113").Value = 1099.5
$ws.Range("K113").Value = 45500.26
$ws.Range("L113").Value = 1099.5
$ws.Range("M113").Value = -43330.26
$ws.Range("N113").Value = -5439.5

$ws.Range("H136").Value = 86964540
$ws.Range("I136").Value = 58831884
$ws.Range("J136").Value = 166673730
$ws.Range("K136").Value = 176495652
$ws.Range("L136").Value = 500021190
$ws.Range("M136").Value = -176493102
$ws.Range("N136").Value = -500026290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 47904.89
$ws.Range("J54").Value = 99015.336
$ws.Range("L54").Value = 99015.336
$ws.Range("N54").Value = -100055.336

$ws.Range("H136").Value = 9123.4375
$ws.Range("I136").Value = 8284
$ws.Range("K136").Value = 24852
$ws.Range("M136").Value = -22302
